$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1656.4546
$ws.Range("I53").Value = 677.75
$ws.Range("J53").Value = 4266.3335
$ws.Range("K53").Value = 677.75
$ws.Range("L53").Value = 4266.3335
$ws.Range("M53").Value = -40.75
$ws.Range("N53").Value = -5540.3335

$ws.Range("H129").Value = 3666.6667
$ws.Range("I129").Value = 3666.6667
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 11000.0001
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -6000.000100000001

$ws.Range("H132").Value = 4669.755
$ws.Range("I132").Value = 2420.2258
$ws.Range("J132").Value = 8543.944
$ws.Range("K132").Value = 7260.6774
$ws.Range("L132").Value = 25631.832
$ws.Range("M132").Value = -4730.6774
$ws.Range("N132").Value = -30691.832

$ws.Range("H135").Value = 65217704
$ws.Range("I135").Value = 22727602
$ws.Range("J135").Value = 1000000000
$ws.Range("K135").Value = 204548418
$ws.Range("L135").Value = 9000000000
$ws.Range("M135").Value = -204545883
$ws.Range("N135").Value = -9000005070

$ws.Range("H138").Value = 1372
$ws.Range("I138").Value = 1372
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4116
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = 1024

$ws.Range("H141").Value = 1159.8572
$ws.Range("I141").Value = 1207.0769
$ws.Range("J141").Value = 546
$ws.Range("K141").Value = 3621.2307
$ws.Range("L141").Value = 1638
$ws.Range("M141").Value = 1558.7693
$ws.Range("N141").Value = -11998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 47623190
$ws.Range("I74").Value = 52634896
$ws.Range("K74").Value = 52634896
$ws.Range("M74").Value = -52634022

$ws.Range("H77").Value = 47623190
$ws.Range("I77").Value = 52634896
$ws.Range("K77").Value = 263174480
$ws.Range("M77").Value = -263170112

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1833.625
$ws.Range("I105").Value = 1879.1666
$ws.Range("J105").Value = 1697
$ws.Range("K105").Value = 1879.1666
$ws.Range("L105").Value = 1697
$ws.Range("M105").Value = -132.1666
$ws.Range("N105").Value = -5191

$ws.Range("H134").Value = 42918892
$ws.Range("I134").Value = 42918892
$ws.Range("K134").Value = 128756676
$ws.Range("M134").Value = -128754141

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4450
$ws.Range("I122").Value = 4483.3
$ws.Range("K122").Value = 13449.9
$ws.Range("M122").Value = -10999.9

$ws.Range("H132").Value = 47622000
$ws.Range("I132").Value = 52633524
$ws.Range("J132").Value = 12499.5
$ws.Range("K132").Value = 157900572
$ws.Range("L132").Value = 37498.5
$ws.Range("M132").Value = -157898042
$ws.Range("N132").Value = -42558.5

$ws.Range("H134").Value = 5683289.5
$ws.Range("I134").Value = 5683289.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17049868.5
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -17047333.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 694.7857
$ws.Range("I26").Value = 44
$ws.Range("J26").Value = 955.1
$ws.Range("K26").Value = 132
$ws.Range("L26").Value = 2865.3
$ws.Range("M26").Value = 156
$ws.Range("N26").Value = -3441.3

$ws.Range("H34").Value = 624.1667
$ws.Range("I34").Value = 419.2
$ws.Range("J34").Value = 1649
$ws.Range("K34").Value = 1257.6
$ws.Range("L34").Value = 4947
$ws.Range("M34").Value = -1173.6
$ws.Range("N34").Value = -5115

$ws.Range("H60").Value = 5888.9287
$ws.Range("I60").Value = 311.25
$ws.Range("K60").Value = 933.75
$ws.Range("M60").Value = -682.75

$ws.Range("H81").Value = 283888.78
$ws.Range("I81").Value = 306874.88
$ws.Range("J81").Value = 100000
$ws.Range("K81").Value = 920624.64
$ws.Range("L81").Value = 300000
$ws.Range("M81").Value = -919501.64
$ws.Range("N81").Value = -302246

$ws.Range("H84").Value = 283888.78
$ws.Range("I84").Value = 306874.88
$ws.Range("J84").Value = 100000
$ws.Range("K84").Value = 2761873.92
$ws.Range("L84").Value = 900000
$ws.Range("M84").Value = -2756257.92
$ws.Range("N84").Value = -911232

$ws.Range("H93").Value = 9000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 9000
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = 27000
$ws.Range("N93").Value = -30744

$ws.Range("H104").Value = 3620
$ws.Range("J104").Value = 3620
$ws.Range("L104").Value = 10860
$ws.Range("N104").Value = -16102

$ws.Range("H109").Value = 1928.4286
$ws.Range("I109").Value = 533.3333
$ws.Range("J109").Value = 2974.75
$ws.Range("K109").Value = 1599.9999
$ws.Range("L109").Value = 8924.25
$ws.Range("M109").Value = -559.9999
$ws.Range("N109").Value = -11004.25

$ws.Range("H115").Value = 9999.5
$ws.Range("I115").Value = 10000
$ws.Range("K115").Value = 30000
$ws.Range("M115").Value = -28825

$ws.Range("H118").Value = 565
$ws.Range("I118").Value = 565
$ws.Range("K118").Value = 1695
$ws.Range("M118").Value = -452

$ws.Range("H139").Value = 1349.2858
$ws.Range("I139").Value = 1222.3077
$ws.Range("K139").Value = 3666.9231
$ws.Range("M139").Value = 1473.0769

$ws.Range("H140").Value = 1435.4546
$ws.Range("I140").Value = 1435.4546
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4306.3638
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = 873.6361999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 10040.333

$ws.Range("H132").Value = 13891547
$ws.Range("I132").Value = 13891547
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 41674641
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -41672111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 26674490
$ws.Range("I132").Value = 30008426
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 90025278
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -90022748
$ws.Range("N132").Value = -14058.5

$ws.Range("H136").Value = 2899.5
$ws.Range("I136").Value = 2800
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 8400
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -5850
$ws.Range("N136").Value = -14097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16130760
$ws.Range("I132").Value = 16130760
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 48392280
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -48389750
